$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row entirely (row 26); everything below shifts up one row.
$ws.Rows(26).Delete()

# The table now has one row too many at the bottom (old "SC 232" duplicate after
# the shift) - remove it so the sheet ends at row 33.
$ws.Rows(34).Delete()

# --- Apply remaining per-cell value corrections to match the target data ---

$ws.Range("F2").Value = ""
$ws.Range("F5").Value = 17.66
$ws.Range("D6").Value = -14.2
$ws.Range("F6").Value = 16.43
$ws.Range("D8").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("D12").Value = -14.1
$ws.Range("F13").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("D17").Value = -14.7
$ws.Range("D18").Value = -15.2
$ws.Range("D19").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("D23").Value = -13.9
$ws.Range("F24").Value = 16.78

$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = 10
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17

$ws.Range("A28").Value = "SC 105"
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = ""

$ws.Range("A29").Value = "SC 119"
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = -13
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06

$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = ""
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = ""
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = ""
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
